# Convert the complex field " m:'some text'.noStyleText() " into a plain
# literal-text rendering of the token stream: "{" + tokens + "}" using
# w:t runs instead of w:fldChar/w:instrText, while preserving the run
# splitting and the FFC000 coloring of the quoted-literal/query tokens
# (TokenIteratorFieldRewriterSplit behaviour).

$d = $word.ActiveDocument

# The field lives in its own paragraph; remember that paragraph's start
# offset so we can re-find it once the field (and the character
# positions after it) have shifted.
$field = $d.Fields(1)
$fieldParagraphStart = $field.Code.Paragraphs(1).Range.Start

# Delete the whole field (fldChar begin/end + instrText runs) - this
# leaves an empty paragraph behind, ready to receive literal text runs.
# Positions shift once the field is gone, so re-resolve the insertion
# point afterwards instead of reusing a pre-delete offset.
$field.Delete()

$insertionPoint = $d.Range($fieldParagraphStart, $fieldParagraphStart).Paragraphs(1).Range.Start
$target = $d.Range($insertionPoint, $insertionPoint)

$openXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>{</w:t></w:r>
            <w:r><w:t>m</w:t></w:r>
            <w:r><w:t>:</w:t></w:r>
            <w:r><w:rPr><w:color w:val="FFC000"/></w:rPr><w:t>'</w:t></w:r>
            <w:r><w:rPr><w:color w:val="FFC000"/></w:rPr><w:t>some text</w:t></w:r>
            <w:r><w:rPr><w:color w:val="FFC000"/></w:rPr><w:t>'.</w:t></w:r>
            <w:r><w:rPr><w:color w:val="FFC000"/></w:rPr><w:t>noStyleText</w:t></w:r>
            <w:r><w:rPr><w:color w:val="FFC000"/></w:rPr><w:t>()</w:t></w:r>
            <w:r><w:t xml:space="preserve">}</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($openXml)
